$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Registre de tests")

# --- Apply the same cell style already used by D3:D10/F2:F10 (style index 2)
# to E3:E10, which previously had no explicit style. We copy the format
# from D3 (which already carries that style) onto E3:E10 without touching
# the underlying boolean values.
$ws.Range("D3").Copy() | Out-Null
$ws.Range("E3:E10").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Update the "Réussi" (success) boolean flags now that unit testing for
# correlation / correlation square has been completed.
$ws.Range("E2").Value = $true
$ws.Range("E3").Value = $true
$ws.Range("E4").Value = $false
$ws.Range("E5").Value = $true
$ws.Range("E6").Value = $true
$ws.Range("E7").Value = $false
$ws.Range("E8").Value = $true
$ws.Range("E9").Value = $true
$ws.Range("E10").Value = $false

# --- Row heights for the trailing blank rows shrink slightly.
$ws.Rows.Item(17).RowHeight = 12.75
$ws.Rows.Item(18).RowHeight = 12.75
$ws.Rows.Item(19).RowHeight = 12.75

# --- Move the active selection to E10, matching where editing finished.
$ws.Range("E10").Select() | Out-Null
